$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its text representation instead of
# being auto-converted to numbers by Excel (values like "0.999", "1.00",
# "7.11" would otherwise be parsed as numeric).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "69.382.37"
$ws.Range("E2").Value = "  +0.04%  "

$ws.Range("D3").Value = "3.671.96"
$ws.Range("E3").Value = "  -0.48%  "

$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "644.79"
$ws.Range("E5").Value = "  -5.41%  "

$ws.Range("D6").Value = "159.90"
$ws.Range("E6").Value = "  +0.18%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("E8").Value = "  +0.27%  "

$ws.Range("E9").Value = "  -0.62%  "

$ws.Range("D10").Value = "7.11"
$ws.Range("E10").Value = "  -0.18%  "

$ws.Range("D11").Value = "0.443"
$ws.Range("E11").Value = "  +0.53%  "

$ws.Range("E12").Value = "  +0.28%  "

$ws.Range("D13").Value = "4.286.41"

$ws.Range("D14").Value = "32.69"
$ws.Range("E14").Value = "  +0.49%  "

$ws.Range("D15").Value = "3.670.01"
$ws.Range("E15").Value = "  -0.50%  "

$ws.Range("D16").Value = "69.355.48"
$ws.Range("E16").Value = "  +0.02%  "

$ws.Range("E17").Value = "  +0.37%  "

$ws.Range("D18").Value = "16.04"
$ws.Range("E18").Value = "  -0.08%  "

$ws.Range("D19").Value = "6.49"
$ws.Range("E19").Value = "  -0.06%  "

$ws.Range("D20").Value = "466.24"
$ws.Range("E20").Value = "  -0.75%  "

$ws.Range("D21").Value = "9.88"
$ws.Range("E21").Value = "  -0.50%  "

$ws.Range("D22").Value = "0.646"
$ws.Range("E22").Value = "  -1.62%  "

$ws.Range("D23").Value = "79.48"
$ws.Range("E23").Value = "  -0.66%  "

$ws.Range("D24").Value = "3.816.51"
$ws.Range("E24").Value = "  -0.50%  "

$ws.Range("E25").Value = "  +0.18%  "

$ws.Range("D26").Value = "0.0000126"
$ws.Range("E26").Value = "  +2.33%  "

$ws.Range("D27").Value = "10.87"
$ws.Range("E27").Value = "  -0.54%  "

$ws.Range("D28").Value = "9.06"
$ws.Range("E28").Value = "  -0.93%  "

$ws.Range("E29").Value = "  -2.71%  "

$ws.Range("D30").Value = "1.72"
$ws.Range("E30").Value = "  -1.46%  "

$ws.Range("D31").Value = "2.01"
$ws.Range("E31").Value = "  +0.65%  "

$ws.Range("D32").Value = "1.00"
$ws.Range("E32").Value = "  -0.44%  "

$ws.Range("D33").Value = "26.71"
$ws.Range("E33").Value = "  -0.98%  "

$ws.Range("D34").Value = "6.47"
$ws.Range("E34").Value = "  -2.44%  "

$ws.Range("E35").Value = "  +4.25%  "

$ws.Range("D36").Value = "3.660.55"
$ws.Range("E36").Value = "  -0.45%  "

$ws.Range("D37").Value = "8.45"
$ws.Range("E37").Value = "  +1.93%  "

$ws.Range("D39").Value = "5.92"
$ws.Range("E39").Value = "  -5.73%  "

$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D40").Value = "177.81"
$ws.Range("E40").Value = "  +4.34%  "

$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  -0.08%  "

$ws.Range("D42").Value = "0.0904"
$ws.Range("E42").Value = "  -0.45%  "

$ws.Range("E43").Value = "  -3.06%  "

$ws.Range("D44").Value = "0.925"
$ws.Range("E44").Value = "  -1.94%  "

$ws.Range("D45").Value = "46.63"
$ws.Range("E45").Value = "  -2.08%  "

$ws.Range("D46").Value = "2.74"
$ws.Range("E46").Value = "  +1.38%  "

$ws.Range("B47").Value = "FLOKI"
$ws.Range("C47").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D47").Value = "0.000270"
$ws.Range("E47").Value = "  -2.29%  "

$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "27.16"
$ws.Range("E48").Value = "  -4.17%  "

$ws.Range("E49").Value = "  -4.64%  "

$ws.Range("D50").Value = "1.26"
$ws.Range("E50").Value = "  -3.62%  "

$ws.Range("D51").Value = "7.85"
$ws.Range("E51").Value = "  +0.58%  "

